$d = $word.ActiveDocument

# Pull the whole document as Flat-OPC WordOpenXML so we can perform a
# surgical, run-level rewrite (this lets us create new <w:r> runs with
# exact <w:rPr> formatting -- something the Range/Font object alone
# cannot fully express, e.g. eastAsia/cs font names and szCs).
$full = $d.Content.WordOpenXML

# ---------------------------------------------------------------
# Edit 1: "mays plustost il les fault bien tremper en " ->
#   "mays plustost il les fault " + <<corr><del> + "les" + </del></corr> + " " + "bien tremper en "
# ---------------------------------------------------------------
$old1 = '<w:r w:rsidDel="00000000" w:rsidRPr="00000000"><w:rPr><w:color w:val="000000"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">mays plustost il les fault bien tremper en </w:t></w:r>'

$new1 = '<w:r><w:rPr><w:color w:val="000000"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">mays plustost il les fault </w:t></w:r>' + `
        '<w:r><w:rPr><w:rFonts w:ascii="Courier New" w:cs="Courier New" w:eastAsia="Courier New" w:hAnsi="Courier New"/><w:color w:val="a91111"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">&lt;&lt;corr&gt;&lt;del&gt;</w:t></w:r>' + `
        '<w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">les</w:t></w:r>' + `
        '<w:r><w:rPr><w:rFonts w:ascii="Courier New" w:cs="Courier New" w:eastAsia="Courier New" w:hAnsi="Courier New"/><w:color w:val="a91111"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">&lt;/del&gt;&lt;/corr&gt; </w:t></w:r>' + `
        '<w:r><w:rPr><w:color w:val="000000"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">bien tremper en </w:t></w:r>'

if ($full.IndexOf($old1) -lt 0) {
    throw "edit1: anchor text not found"
}
$full = $full.Replace($old1, $new1)

# ---------------------------------------------------------------
# Edit 2: "aprime rectificatam, deinde pulvere composito aspergatur" ->
#   "aprime rectificat" + "u" + "m, deinde pulvere composito aspergatur"
# ---------------------------------------------------------------
$old2 = '<w:r w:rsidDel="00000000" w:rsidRPr="00000000"><w:rPr><w:color w:val="000000"/><w:rtl w:val="0"/></w:rPr><w:t>aprime rectificatam, deinde pulvere composito aspergatur</w:t></w:r>'

$new2 = '<w:r><w:rPr><w:color w:val="000000"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">aprime rectificat</w:t></w:r>' + `
        '<w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">u</w:t></w:r>' + `
        '<w:r><w:rPr><w:color w:val="000000"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">m, deinde pulvere composito aspergatur</w:t></w:r>'

if ($full.IndexOf($old2) -lt 0) {
    throw "edit2: anchor text not found"
}
$full = $full.Replace($old2, $new2)

$null = $d.Content.InsertXML($full)
